$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=2; D='26.201.45'; E='  +0.25%  ' },
    @{ Row=3; D='1.651.88'; E='  -0.30%  ' },
    @{ Row=4; D='1.006'; E='  +0.53%  ' },
    @{ Row=5; D='218.06'; E='  +1.39%  ' },
    @{ Row=6; D='0.5191'; E='  -0.88%  ' },
    @{ Row=7; D='1.006'; E='  +0.46%  ' },
    @{ Row=8; D='0.2643'; E='  +0.65%  ' },
    @{ Row=9; D='0.06294'; E='  -1.58%  ' },
    @{ Row=10; D='21.25'; E='  +1.90%  ' },
    @{ Row=11; D='0.07733'; E='  -0.27%  ' },
    @{ Row=12; D='1.675.16'; E='  +1.08%  ' },
    @{ Row=13; D='4.415'; E='  -0.82%  ' },
    @{ Row=14; D='0.5438'; E='  -1.52%  ' },
    @{ Row=15; D='0.0₅8195'; E='  -1.02%  ' },
    @{ Row=16; D='64.60'; E='  -0.89%  ' },
    @{ Row=17; D='26.227.60'; E='  +0.27%  ' },
    @{ Row=18; E='  +0.40%  ' },
    @{ Row=19; D='4.673'; E='  -1.93%  ' },
    @{ Row=20; D='190.47'; E='  +0.10%  ' },
    @{ Row=21; D='10.17'; E='  -1.45%  ' },
    @{ Row=22; D='6.180'; E='  -3.01%  ' },
    @{ Row=23; E='  +0.56%  ' },
    @{ Row=24; D='138.58'; E='  -3.23%  ' },
    @{ Row=25; D='0.1243'; E='  -0.58%  ' },
    @{ Row=26; D='7.271'; E='  -1.82%  ' },
    @{ Row=27; D='16.03'; E='  +0.30%  ' },
    @{ Row=28; D='1.415'; E='  -1.45%  ' },
    @{ Row=29; D='0.06047'; E='  -1.35%  ' },
    @{ Row=30; D='1.283'; E='  +1.51%  ' },
    @{ Row=31; D='3.545'; E='  +1.21%  ' },
    @{ Row=32; D='3.349'; E='  -2.08%  ' },
    @{ Row=33; D='1.650'; E='  -0.72%  ' },
    @{ Row=34; D='0.9832'; E='  -1.60%  ' },
    @{ Row=35; D='2.413'; E='  +0.52%  ' },
    @{ Row=36; D='2.775'; E='  +0.42%  ' },
    @{ Row=37; D='0.5928'; E='  +4.75%  ' },
    @{ Row=38; D='0.01595'; E='  -0.68%  ' },
    @{ Row=39; D='5.953'; E='  +0.85%  ' },
    @{ Row=40; B='TrustWalletToken'; C='https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; D='0.8619'; E='  +0.90%  ' },
    @{ Row=41; B='Maker'; C='https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'; D='1.059.42'; E='  +2.68%  ' },
    @{ Row=42; E='  +0.25%  ' },
    @{ Row=43; D='99.88'; E='  +0.11%  ' },
    @{ Row=44; D='1.797.34'; E='  -0.41%  ' },
    @{ Row=45; D='0.0₈108'; E='  +1.50%  ' },
    @{ Row=46; D='57.22'; E='  +2.11%  ' },
    @{ Row=47; D='1.003'; E='  -0.05%  ' },
    @{ Row=48; D='8.059'; E='  -0.70%  ' },
    @{ Row=49; D='0.05182'; E='  +0.38%  ' },
    @{ Row=50; D='1.467'; E='  +4.90%  ' },
    @{ Row=51; E='  +0.42%  ' }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    if ($r.ContainsKey("B")) { $ws.Cells.Item($rowNum, 2).Value = $r.B }
    if ($r.ContainsKey("C")) { $ws.Cells.Item($rowNum, 3).Value = $r.C }
    if ($r.ContainsKey("D")) {
        $cell = $ws.Cells.Item($rowNum, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $r.D
        $cell.ClearFormats()
    }
    if ($r.ContainsKey("E")) { $ws.Cells.Item($rowNum, 5).Value = $r.E }
}
